$wb = $excel.ActiveWorkbook

# Sheet: 展览 (sheet1)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 147
$ws.Range("F3").Value = 922
$ws.Range("F6").Value = 336
$ws.Range("F8").Value = 12394
$ws.Range("F16").Value = 226
$ws.Range("F17").Value = 281
$ws.Range("F18").Value = 790
$ws.Range("F21").Value = 2935
$ws.Range("F23").Value = 4289
$ws.Range("F24").Value = 1144
$ws.Range("F29").Value = 1082
$ws.Range("F31").Value = 112
$ws.Range("F32").Value = 279
$ws.Range("F36").Value = 20
$ws.Range("F37").Value = 4476
$ws.Range("F39").Value = 4588
$ws.Range("F40").Value = 5570
$ws.Range("F45").Value = 364
$ws.Range("F46").Value = 83
$ws.Range("F48").Value = 4114
$ws.Range("F49").Value = 140

# Sheet: 演出 (sheet2)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F12").Value = 1046

# Sheet: 本地生活 (sheet3)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 96
$ws.Range("F5").Value = 13

# Sheet: 全部类型 (sheet4)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 147
$ws.Range("F4").Value = 922
$ws.Range("F7").Value = 336
$ws.Range("F9").Value = 12394
$ws.Range("F14").Value = 226
$ws.Range("F15").Value = 281
$ws.Range("F16").Value = 790
$ws.Range("F19").Value = 2935
$ws.Range("F21").Value = 4289
$ws.Range("F22").Value = 4289
$ws.Range("F23").Value = 1144
$ws.Range("F31").Value = 1082
$ws.Range("F33").Value = 112
$ws.Range("F35").Value = 279
$ws.Range("F39").Value = 4588
$ws.Range("F45").Value = 83
$ws.Range("F46").Value = 4114
